$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue "D2" '68.843.07'
Set-TextValue "E2" '  -0.18%  '

Set-TextValue "D3" '3.856.54'
Set-TextValue "E3" '  +2.98%  '

Set-TextValue "E4" '  +0.08%  '

Set-TextValue "D5" '601.97'
Set-TextValue "E5" '  -0.10%  '

Set-TextValue "D6" '162.52'
Set-TextValue "E6" '  -2.91%  '

Set-TextValue "D7" '3.856.09'
Set-TextValue "E7" '  +3.03%  '

Set-TextValue "E8" '  +0.08%  '

Set-TextValue "D9" '0.530'
Set-TextValue "E9" '  -1.58%  '

Set-TextValue "D10" '0.167'
Set-TextValue "E10" '  -1.32%  '

Set-TextValue "D11" '6.30'
Set-TextValue "E11" '  -2.76%  '

Set-TextValue "E12" '  -0.25%  '

Set-TextValue "D13" '36.81'
Set-TextValue "E13" '  -2.96%  '

Set-TextValue "E14" '  -2.11%  '

Set-TextValue "D15" '4.508.88'
Set-TextValue "E15" '  +3.10%  '

Set-TextValue "D16" '3.869.45'
Set-TextValue "E16" '  +3.12%  '

Set-TextValue "D17" '69.013.22'
Set-TextValue "E17" '  +0.03%  '

Set-TextValue "E18" '  +2.57%  '

Set-TextValue "E19" '  -0.22%  '

Set-TextValue "B20" 'Chainlink'
Set-TextValue "C20" 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextValue "D20" '17.10'
Set-TextValue "E20" '  -0.86%  '

Set-TextValue "B21" 'Uniswap'
Set-TextValue "C21" 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
Set-TextValue "D21" '11.35'
Set-TextValue "E21" '  +4.46%  '

Set-TextValue "D22" '483.71'
Set-TextValue "E22" '  -1.73%  '

Set-TextValue "D23" '0.720'
Set-TextValue "E23" '  -0.76%  '

Set-TextValue "E24" '  +4.32%  '

Set-TextValue "D25" '83.89'
Set-TextValue "E25" '  -1.02%  '

Set-TextValue "D26" '2.25'
Set-TextValue "E26" '  -2.57%  '

Set-TextValue "D27" '12.12'
Set-TextValue "E27" '  -1.90%  '

Set-TextValue "D28" '10.03'
Set-TextValue "E28" '  -0.83%  '

Set-TextValue "E29" '  -0.12%  '

Set-TextValue "D30" '2.97'
Set-TextValue "E30" '  -0.63%  '

Set-TextValue "E31" '  -0.70%  '

Set-TextValue "D32" '4.011.03'
Set-TextValue "E32" '  +3.08%  '

Set-TextValue "E33" '  -3.74%  '

Set-TextValue "D34" '32.20'
Set-TextValue "E34" '  +2.07%  '

Set-TextValue "D35" '3.807.90'
Set-TextValue "E35" '  +3.47%  '

Set-TextValue "E36" '  -1.57%  '

Set-TextValue "D37" '1.03'
Set-TextValue "E37" '  +1.29%  '

Set-TextValue "E38" '  +4.02%  '

Set-TextValue "E39" '  +0.00%  '

Set-TextValue "E40" '  +0.03%  '

Set-TextValue "E41" '  -1.86%  '

Set-TextValue "D42" '442.49'
Set-TextValue "E42" '  +2.79%  '

Set-TextValue "D43" '2.98'
Set-TextValue "E43" '  +0.16%  '

Set-TextValue "D44" '48.46'
Set-TextValue "E44" '  -0.35%  '

Set-TextValue "E45" '  -0.78%  '

Set-TextValue "E46" '  +0.00%  '

Set-TextValue "D47" '8.39'
Set-TextValue "E47" '  -0.98%  '

Set-TextValue "D48" '26.71'
Set-TextValue "E48" '  +13.15%  '

Set-TextValue "D49" '143.24'
Set-TextValue "E49" '  +1.54%  '

Set-TextValue "D50" '2.831.01'
Set-TextValue "E50" '  +1.72%  '

Set-TextValue "D51" '0.0357'
Set-TextValue "E51" '  +1.53%  '
